$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("data_RAM code_FLASH").Name = "data_RAM-code_FLASH"
$wb.Worksheets.Item("data_CCM code_FLASH").Name = "data_CCM-code_FLASH"
$wb.Worksheets.Item("data_RAM code_CCM").Name = "data_RAM-code_CCM"
$wb.Worksheets.Item("data_CCM code_CCM").Name = "data_CCM-code_CCM"
